# The sheet "Hortaliza, Femacal de La Calera - Cilantro" gets a new weekly
# price-report row inserted at row 188 (pushing the former rows 188..249
# down to 189..250, growing the used range from A1:R249 to A1:R250).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 188; everything below shifts down one.
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row with the new report entry.
$ws.Range("A188").Value = 3
$ws.Range("B188").Value = "Femacal de La Calera"
$ws.Range("C188").Value = "Coquimbo"
$ws.Range("D188").Value = 44524
$ws.Range("E188").Value = 5
$ws.Range("F188").Value = 100112040
$ws.Range("G188").Value = "Cilantro"
$ws.Range("H188").Value = "Sin especificar"
$ws.Range("I188").Value = "Primera"
$ws.Range("J188").Value = 120
$ws.Range("K188").Value = 3000
$ws.Range("L188").Value = 3000
$ws.Range("M188").Value = 3000
$ws.Range("N188").Value = "$/docena de atados (3 kilos)"
$ws.Range("O188").Value = "Provincia de Quillota"
$ws.Range("P188").Value = 1000
$ws.Range("Q188").Value = 3
$ws.Range("R188").Value = "Hortaliza"
